# Auto-generated Excel COM-interop edit script
# Applies the crypto price/volume refresh described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without letting Excel
# silently reinterpret numeric-looking strings (e.g. '1.00') as numbers,
# which would drop the trailing zero / formatting we need to preserve.
function Set-CellText {
    param($Cell, [string]$Text)
    if ($Text -match '^-?[0-9]+(\.[0-9]+)?$') {
        $Cell.NumberFormat = '@'
    }
    $Cell.Value = $Text
}

# Row 2
Set-CellText $ws.Cells.Item(2, 4) '36.624.32'
Set-CellText $ws.Cells.Item(2, 5) '  +0.35%  '

# Row 3
Set-CellText $ws.Cells.Item(3, 4) '2.033.14'
Set-CellText $ws.Cells.Item(3, 5) '  +3.68%  '

# Row 4
Set-CellText $ws.Cells.Item(4, 5) '  -0.06%  '

# Row 5
Set-CellText $ws.Cells.Item(5, 4) '232.98'
Set-CellText $ws.Cells.Item(5, 5) '  -6.92%  '

# Row 6
Set-CellText $ws.Cells.Item(6, 5) '  -0.46%  '

# Row 7
Set-CellText $ws.Cells.Item(7, 5) '  +0.00%  '

# Row 8
Set-CellText $ws.Cells.Item(8, 4) '55.42'
Set-CellText $ws.Cells.Item(8, 5) '  +2.78%  '

# Row 9
Set-CellText $ws.Cells.Item(9, 4) '0.372'
Set-CellText $ws.Cells.Item(9, 5) '  +0.90%  '

# Row 10
Set-CellText $ws.Cells.Item(10, 4) '57.16'
Set-CellText $ws.Cells.Item(10, 5) '  +3.30%  '

# Row 11
Set-CellText $ws.Cells.Item(11, 5) '  -0.17%  '

# Row 12
Set-CellText $ws.Cells.Item(12, 5) '  +0.21%  '

# Row 13
Set-CellText $ws.Cells.Item(13, 4) '2.330.82'
Set-CellText $ws.Cells.Item(13, 5) '  +3.67%  '

# Row 14
Set-CellText $ws.Cells.Item(14, 4) '14.36'
Set-CellText $ws.Cells.Item(14, 5) '  +3.09%  '

# Row 15
Set-CellText $ws.Cells.Item(15, 4) '20.25'
Set-CellText $ws.Cells.Item(15, 5) '  -3.89%  '

# Row 16
Set-CellText $ws.Cells.Item(16, 4) '0.764'
Set-CellText $ws.Cells.Item(16, 5) '  +0.72%  '

# Row 17
Set-CellText $ws.Cells.Item(17, 5) '  +1.24%  '

# Row 18
Set-CellText $ws.Cells.Item(18, 4) '2.028.31'
Set-CellText $ws.Cells.Item(18, 5) '  +3.51%  '

# Row 19
Set-CellText $ws.Cells.Item(19, 4) '36.782.94'
Set-CellText $ws.Cells.Item(19, 5) '  +1.11%  '

# Row 20
Set-CellText $ws.Cells.Item(20, 4) '67.55'
Set-CellText $ws.Cells.Item(20, 5) '  -1.27%  '

# Row 21
Set-CellText $ws.Cells.Item(21, 4) '5.54'
Set-CellText $ws.Cells.Item(21, 5) '  +11.31%  '

# Row 22
Set-CellText $ws.Cells.Item(22, 4) '0.0₃0798'
Set-CellText $ws.Cells.Item(22, 5) '  -1.97%  '

# Row 23
Set-CellText $ws.Cells.Item(23, 4) '221.07'
Set-CellText $ws.Cells.Item(23, 5) '  -3.73%  '

# Row 24
Set-CellText $ws.Cells.Item(24, 5) '  +0.03%  '

# Row 25
Set-CellText $ws.Cells.Item(25, 4) '2.39'
Set-CellText $ws.Cells.Item(25, 5) '  +1.06%  '

# Row 26
Set-CellText $ws.Cells.Item(26, 5) '  -4.14%  '

# Row 27
Set-CellText $ws.Cells.Item(27, 4) '162.89'
Set-CellText $ws.Cells.Item(27, 5) '  +0.03%  '

# Row 28
Set-CellText $ws.Cells.Item(28, 5) '  -0.17%  '

# Row 29
Set-CellText $ws.Cells.Item(29, 5) '  +6.99%  '

# Row 30
Set-CellText $ws.Cells.Item(30, 4) '18.93'
Set-CellText $ws.Cells.Item(30, 5) '  -0.34%  '

# Row 31
Set-CellText $ws.Cells.Item(31, 5) '  +3.43%  '

# Row 32
Set-CellText $ws.Cells.Item(32, 5) '  +0.30%  '

# Row 33
Set-CellText $ws.Cells.Item(33, 5) '  -1.38%  '

# Row 34
Set-CellText $ws.Cells.Item(34, 4) '0.0603'
Set-CellText $ws.Cells.Item(34, 5) '  -1.87%  '

# Row 35
Set-CellText $ws.Cells.Item(35, 5) '  +6.01%  '

# Row 36
Set-CellText $ws.Cells.Item(36, 4) '4.28'
Set-CellText $ws.Cells.Item(36, 5) '  +0.40%  '

# Row 37
Set-CellText $ws.Cells.Item(37, 4) '1.00'
Set-CellText $ws.Cells.Item(37, 5) '  -0.21%  '

# Row 38
Set-CellText $ws.Cells.Item(38, 2) 'THORChain'
Set-CellText $ws.Cells.Item(38, 3) 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-CellText $ws.Cells.Item(38, 4) '5.86'
Set-CellText $ws.Cells.Item(38, 5) '  +11.77%  '

# Row 39
Set-CellText $ws.Cells.Item(39, 4) '3.30'
Set-CellText $ws.Cells.Item(39, 5) '  -0.39%  '

# Row 40
Set-CellText $ws.Cells.Item(40, 2) 'WEMIXToken'
Set-CellText $ws.Cells.Item(40, 3) 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-CellText $ws.Cells.Item(40, 4) '1.76'
Set-CellText $ws.Cells.Item(40, 5) '  -2.62%  '

# Row 41
Set-CellText $ws.Cells.Item(41, 5) '  -2.08%  '

# Row 42
Set-CellText $ws.Cells.Item(42, 4) '1.474.43'
Set-CellText $ws.Cells.Item(42, 5) '  +2.70%  '

# Row 43
Set-CellText $ws.Cells.Item(43, 4) '0.0931'
Set-CellText $ws.Cells.Item(43, 5) '  +3.77%  '

# Row 44
Set-CellText $ws.Cells.Item(44, 4) '92.95'
Set-CellText $ws.Cells.Item(44, 5) '  +7.15%  '

# Row 45
Set-CellText $ws.Cells.Item(45, 2) 'TrustWalletToken'
Set-CellText $ws.Cells.Item(45, 3) 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-CellText $ws.Cells.Item(45, 4) '1.13'
Set-CellText $ws.Cells.Item(45, 5) '  -1.13%  '

# Row 46
Set-CellText $ws.Cells.Item(46, 2) 'FTXToken'
Set-CellText $ws.Cells.Item(46, 3) 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-CellText $ws.Cells.Item(46, 4) '4.15'
Set-CellText $ws.Cells.Item(46, 5) '  +37.91%  '

# Row 47
Set-CellText $ws.Cells.Item(47, 4) '0.0204'
Set-CellText $ws.Cells.Item(47, 5) '  +0.10%  '

# Row 48
Set-CellText $ws.Cells.Item(48, 4) '15.70'
Set-CellText $ws.Cells.Item(48, 5) '  +2.92%  '

# Row 49
Set-CellText $ws.Cells.Item(49, 4) '1.01'
Set-CellText $ws.Cells.Item(49, 5) '  +1.00%  '

# Row 50
Set-CellText $ws.Cells.Item(50, 5) '  +1.62%  '

# Row 51
Set-CellText $ws.Cells.Item(51, 4) '6.89'
Set-CellText $ws.Cells.Item(51, 5) '  +2.30%  '

Write-Host "Applied crypto price/volume refresh."
